$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.99
$wsSummary.Range("B4").Value = -0.02
$wsSummary.Range("B6").Value = 106
$wsSummary.Range("B8").Value = 52
$wsSummary.Range("B9").Value = 37.74

# ---- Strategy Status sheet (MarketMaking row) ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.98999999999999
$wsStatus.Range("D4").Value = 106
$wsStatus.Range("E4").Value = -0.02
$wsStatus.Range("F4").Value = -0.01
$wsStatus.Range("G4").Value = 37.74

# ---- New trade row (Trade #106) appended to "All Trades" and "MarketMaking" sheets ----
$newRow = @(106, "2026-02-17", "15:59:00", "MarketMaking", "UP", 0.03, 0.02, "CLOSED", -33.3333, -0.01, 99.98999999999999, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.15)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 107
    $ws.Cells.Item($row, 1).Value = $newRow[0]
    # Date/time look-alike strings must stay text (matches the rest of the column)
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $newRow[1]
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $newRow[2]
    $ws.Cells.Item($row, 4).Value = $newRow[3]
    $ws.Cells.Item($row, 5).Value = $newRow[4]
    $ws.Cells.Item($row, 6).Value = $newRow[5]
    $ws.Cells.Item($row, 7).Value = $newRow[6]
    $ws.Cells.Item($row, 8).Value = $newRow[7]
    $ws.Cells.Item($row, 9).Value = $newRow[8]
    $ws.Cells.Item($row, 10).Value = $newRow[9]
    $ws.Cells.Item($row, 11).Value = $newRow[10]
    $ws.Cells.Item($row, 12).Value = $newRow[11]
    $ws.Cells.Item($row, 13).Value = $newRow[12]
    $ws.Cells.Item($row, 14).Value = $newRow[13]
    $ws.Cells.Item($row, 15).Value = $newRow[14]
    $ws.Cells.Item($row, 16).Value = $newRow[15]
    $ws.Cells.Item($row, 17).Value = $newRow[16]
}
